$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q4" sheet right after "总计" (i.e. before the
#    current "2022-Q3" sheet, which is Worksheets.Item(2)).
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item(2)
$newWs = $wb.Worksheets.Add($q3)
$newWs.Name = "2022-Q4"

# Header row
$newWs.Range("B1").Value = "基金代码"
$newWs.Range("C1").Value = "基金名称"
$newWs.Range("D1").Value = "基金规模"
$newWs.Range("E1").Value = "股票总仓位"
$newWs.Range("F1").Value = "仓位占比"
$newWs.Range("G1").Value = "持有市值(亿元)"
$newWs.Range("H1").Value = "仓位排名"

# Columns that must stay text (numeric-looking values with meaningful
# leading/trailing zeros, e.g. fund codes "000043" and sizes "0.10").
$newWs.Range("B2:B7").NumberFormat = "@"
$newWs.Range("D2:G7").NumberFormat = "@"

$rows = @(
  @("000043", "嘉实美国成长股票（QDII）人民币",       "12.69", "92.23", "1.58", "0.2005", 7),
  @("000044", "嘉实美国成长股票（QDII）美元现汇",       "12.69", "92.23", "1.58", "0.2005", 7),
  @("000369", "广发全球医疗保健（QDII）人民币A",        "3.16",  "80.87", "3.82", "0.1207", 5),
  @("000370", "广发全球医疗保健（QDII）美元A",          "3.16",  "80.87", "3.82", "0.1207", 5),
  @("016280", "广发全球医疗保健（QDII）人民币C",        "0.10",  "80.87", "3.82", "0.0038", 5),
  @("016281", "广发全球医疗保健（QDII）美元C",          "0.10",  "80.87", "3.82", "0.0038", 5)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $newWs.Range("A$r").Value = $i
    $newWs.Range("B$r").Value = $rows[$i][0]
    $newWs.Range("C$r").Value = $rows[$i][1]
    $newWs.Range("D$r").Value = $rows[$i][2]
    $newWs.Range("E$r").Value = $rows[$i][3]
    $newWs.Range("F$r").Value = $rows[$i][4]
    $newWs.Range("G$r").Value = $rows[$i][5]
    $newWs.Range("H$r").Value = $rows[$i][6]
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: insert a new row for 2022-Q4 at the
#    top of the data block and renumber the running index column.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$total.Rows.Item(2).Insert()

# Pick up the index-column formatting (style) from the row underneath before
# overwriting any values, then write the new totals row.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 6
$total.Range("D2").Value = 0.65
$total.Range("B2:D2").ClearFormats()

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
$total.Range("A8").Value = 6
